$d = $word.ActiveDocument

# Locate the "Clean up this list (again)." list item.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Clean up this list (again).*") {
        $target = $p
        break
    }
}

# The trailing _GoBack bookmark currently sits right at the end of this
# paragraph (immediately before its paragraph mark); remove it so it can
# be recreated after the newly inserted paragraph.
$hadBookmark = $d.Bookmarks.Exists("_GoBack")
if ($hadBookmark) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Insert a new list paragraph right after it, inheriting the same
# ListParagraph style / numbering, with the new todo text.
# A trailing placeholder character is appended temporarily: adding a
# bookmark at a position that lands exactly on "end of paragraph" (i.e.
# immediately before the paragraph mark) doesn't land correctly, so the
# bookmark is added just before the placeholder (a safe, non-edge
# position) and the placeholder is deleted afterwards, leaving the
# bookmark correctly collapsed at the end of the real text.
$target.Range.InsertParagraphAfter()
$newPara = $target.Next()
$newPara.Range.Text = "Forward declaration header.~"

if ($hadBookmark) {
    $bmPos = $newPara.Range.End - 2
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

$placeholder = $d.Range($newPara.Range.End - 2, $newPara.Range.End - 1)
$placeholder.Delete()
